# [GQA] Atualização do Registro de NC.
# Row 3 (NC #1) escalates: Criticidade Baixa -> Média, Status Relatada -> Reescalonada.
# A new "Alta" entry is added to the Criticidade legend (purple fill), mirroring
# the existing Baixa/Média/Solucionada legend cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Criticidade (C3): Baixa -> Média ---
# Reuse the existing "Média" legend formatting (K6) so the blue fill / font /
# border match the other Média cells in the table (C4:C6 stay "Baixa").
$ws.Range("K6").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Média"

# --- Status (F3): Relatada -> Reescalonada ---
# Reuse the existing "Reescalonada" legend formatting (I6) for the red fill.
$ws.Range("I6").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "Reescalonada"

$excel.CutCopyMode = $false

# --- Criticidade legend: add "Alta" entry under "Solucionada" (K7) ---
# Start from the same look as the other legend swatches (font/border/alignment)
# then recolor it purple for the new "Alta" level.
$ws.Range("K6").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("K7").Value = "Alta"
$ws.Range("K7").Interior.Color = 10498160

$excel.CutCopyMode = $false

# --- Selection / view: move active cell to A7, no pinned top-left cell ---
[void]$ws.Range("A7").Select()
